$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that are missing a checkmark ("ü", rendered via Wingdings as a tick)
# in the "Tugas" (task) completion grid. Fill them in and match the
# Wingdings checkmark formatting used by every other completed cell
# (e.g. C3/D3/E3/F3).
$targets = @("E4", "E14", "E19", "C34")

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $cell.Value = "ü"
    $cell.Font.Name = "Wingdings"
    $cell.Font.Size = 12
}
